# Edit script: add new "ocr(image,saveVar)" command to the image command list,
# add a new "tn.5250" command category (close/open/saveText/typeKeys/updateScreenFields),
# and rename colorbit's first parameter from "source" to "image".
# All of the underlying command lists live on the hidden "#system" worksheet and are
# surfaced through named ranges used by data-validation dropdowns on the "Scenario" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1) Insert a new column before column Z for the new "tn.5250" command list.
#    This shifts the existing web/webalert/webcookie/ws/ws.async/xml columns
#    (Z:AE) one column to the right (AA:AF).
# ---------------------------------------------------------------------------
$ws.Range("Z1").EntireColumn.Insert()

$ws.Range("Z1").Value2 = "tn.5250"
$ws.Range("Z2").Value2 = "close(profile)"
$ws.Range("Z3").Value2 = "open(profile)"
$ws.Range("Z4").Value2 = "saveText(profile,var)"
$ws.Range("Z5").Value2 = "typeKeys(profile,keystrokes)"
$ws.Range("Z6").Value2 = "updateScreenFields(profile)"

# ---------------------------------------------------------------------------
# 2) "image" command list (column K): rename colorbit's param, and insert the
#    new "ocr(image,saveVar)" command in its alphabetically-sorted position
#    (between "crop(...)" at K5 and "resize(...)" at K6).
# ---------------------------------------------------------------------------
$ws.Range("K2").Value2 = "colorbit(image,bit,saveTo)"

$ws.Range("K8").Value2 = $ws.Range("K7").Value2
$ws.Range("K7").Value2 = $ws.Range("K6").Value2
$ws.Range("K6").Value2 = "ocr(image,saveVar)"

# ---------------------------------------------------------------------------
# 3) "target" list (column A): insert the new "tn.5250" entry in its
#    alphabetically-sorted position (between "step" at A25 and "web" at A26).
# ---------------------------------------------------------------------------
$ws.Range("A32").Value2 = $ws.Range("A31").Value2
$ws.Range("A31").Value2 = $ws.Range("A30").Value2
$ws.Range("A30").Value2 = $ws.Range("A29").Value2
$ws.Range("A29").Value2 = $ws.Range("A28").Value2
$ws.Range("A28").Value2 = $ws.Range("A27").Value2
$ws.Range("A27").Value2 = $ws.Range("A26").Value2
$ws.Range("A26").Value2 = "tn.5250"

# ---------------------------------------------------------------------------
# 4) Update the defined names (named ranges) to reflect the new layout.
# ---------------------------------------------------------------------------
$wb.Names.Item("image").RefersTo      = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo     = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo        = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo   = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo  = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo         = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo   = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo        = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")

# ---------------------------------------------------------------------------
# 5) Restore the originally-active sheet ("Scenario").
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Scenario").Activate()
